# Preplanning of sprint 3
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Sprint Backlog: close out Sprint 2 items, add Sprint 3 backlog
# ---------------------------------------------------------------
$sprint = $wb.Worksheets.Item("Sprint Backlog")
$sprint.Activate()

# Mark the remaining sprint-2 stories as done
$sprint.Range("L5").Value = "done"
$sprint.Range("L6").Value = "done"
$sprint.Range("L7").Value = "done"
$sprint.Range("L8").Value = "done"
$sprint.Range("L9").Value = "done"
$sprint.Range("L10").Value = "done"
$sprint.Range("L11").Value = "done"
$sprint.Range("L12").Value = "done"

# Actual effort recorded for the "Medication Data Model" rework
$sprint.Range("K6").Value = 5

# New Sprint 3 backlog items
$sprint.Range("A13").Value = 4.1
$sprint.Range("B13").Value = 3
$sprint.Range("C13").Value = "Unit Test's"
$sprint.Range("D13").Value = "Create for the logical classes unit test's for testing"
$sprint.Range("E13").Value = "Test"
$sprint.Range("I13").Value = 10
$sprint.Range("L13").Value = "work in progress"

$sprint.Range("A14").Value = 1.5
$sprint.Range("B14").Value = 3
$sprint.Range("C14").Value = "Medication Photo"
$sprint.Range("D14").Value = "Implement the photo of the medication"
$sprint.Range("E14").Value = "UI"
$sprint.Range("I14").Value = 10
$sprint.Range("L14").Value = "work in progress"

$sprint.Range("A15").Value = 2.5
$sprint.Range("B15").Value = 3
$sprint.Range("C15").Value = "Prescription historization"
$sprint.Range("D15").Value = "historization the Insert, editing, and deleting of the Prescription "
$sprint.Range("E15").Value = "Modell, Database"
$sprint.Range("I15").Value = 10
$sprint.Range("L15").Value = "work in progress"

$sprint.Range("A16").Value = 5.1
$sprint.Range("B16").Value = 3
$sprint.Range("C16").Value = "Medication extern Information"
$sprint.Range("D16").Value = "Insert a Compendium link for additional information"
$sprint.Range("E16").Value = "UI, Controller"
$sprint.Range("I16").Value = 10
$sprint.Range("L16").Value = "work in progress"

$sprint.Range("A17").Value = 1.4
$sprint.Range("B17").Value = 3
$sprint.Range("C17").Value = "Login"
$sprint.Range("D17").Value = "Create a login page for the patient"
$sprint.Range("E17").Value = "UI, Controller"
$sprint.Range("I17").Value = 5
$sprint.Range("L17").Value = "work in progress"

$sprint.Range("A18").Value = 6.1
$sprint.Range("B18").Value = 3
$sprint.Range("C18").Value = "Alert Popup"
$sprint.Range("D18").Value = "Create a alert popup for the daily medication with applied funktion"
$sprint.Range("E18").Value = "UI, Controller"
$sprint.Range("I18").Value = 10
$sprint.Range("L18").Value = "work in progress"

$sprint.Range("D22").Select() | Out-Null

# ---------------------------------------------------------------
# 2) Product Backlog: mark preplanning status per story
# ---------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")
$backlog.Activate()

$backlog.Range("H2").Value = "OK"
$backlog.Range("H3").Value = "OK"
$backlog.Range("H4").Value = "OK"
$backlog.Range("H5").Value = "OK"
$backlog.Range("H6").Value = "sprint 3"
$backlog.Range("H7").Value = "sprint 3"
$backlog.Range("H8").Value = "not used"
$backlog.Range("H9").Value = "sprint 3"
$backlog.Range("H10").Value = "OK"
$backlog.Range("H11").Value = "sprint 3"
$backlog.Range("H12").Value = "sprint 3"

$backlog.Range("H21").Select() | Out-Null

# ---------------------------------------------------------------
# 3) Leave Sprint Backlog as the active sheet/tab
# ---------------------------------------------------------------
$sprint.Activate()
